$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells that look numeric stay as text (matches original inlineStr cell type)
$textCells = @("D5", "D6", "D9", "D14", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D31", "D38", "D39", "D40", "D43", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "64.457.31"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.142.48"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "608.79"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "144.21"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D8").Value = "3.140.72"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").Value = "35.45"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "3.656.29"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "64.416.52"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "3.141.67"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "478.71"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "14.82"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "0.717"
$ws.Range("D23").Value = "7.75"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "85.52"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").Value = "13.44"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "8.47"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +8.42%  "
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "2.06"
$ws.Range("E31").Value = "  -5.19%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "0.0₃0767"
$ws.Range("E37").Value = "  +5.56%  "
$ws.Range("D38").Value = "52.51"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").Value = "447.14"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "8.27"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "2.885.03"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "26.29"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "119.69"
$ws.Range("E51").Value = "  +0.37%  "
